$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Item numbers)
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 5
$ws.Range("A5").Value = 6
$ws.Range("A6").Value = 7
$ws.Range("A7").Value = 8
$ws.Range("A8").Value = 9
$ws.Range("A9").Value = 10
$ws.Range("A10").Value = 11
$ws.Range("A11").Value = 13
$ws.Range("A12").Value = 14
$ws.Range("A13").Value = 16
$ws.Range("A14").Value = 17
$ws.Range("A15").Value = 18

# Column B
$ws.Range("B2").Value = "HEPARINA SODICA SUBCUT 5000UI"
$ws.Range("B3").Value = "PROMETAZINA 25MG, CLORIDRATO"
$ws.Range("B4").Value = "HALOPERIDOL 5MG"
$ws.Range("B5").Value = "CLORPROMAZINA 40MG/ML SOL ORAL"
$ws.Range("B6").Value = "HALOPERIDOL 2MG/ML SOL ORAL"
$ws.Range("B7").Value = "CLORPROMAZINA 25MG"
$ws.Range("B8").Value = "CODEINA 30MG"
$ws.Range("B9").Value = "IMIPRAMINA 25MG"
$ws.Range("B10").Value = "RISPERIDONA 3MG"
$ws.Range("B11").Value = "RISPERIDONA 1MG"
$ws.Range("B12").Value = "LEVOMEPROMAZINA 4% GOTAS"
$ws.Range("B13").Value = "LIDOCAINA 2% C/ VASO CONSTRITO"
$ws.Range("B14").Value = "NITRATO DE CERIO +SULFADIAZINA"
$ws.Range("B15").Value = "COLAGENASE+CLORAFENICOL POMADA 30g"

# Column C
$ws.Range("C2").Value = "sodica;heparina"
$ws.Range("C3").Value = "PROMETAZINA"
$ws.Range("C4").Value = "HALOPERIDOL"
$ws.Range("C5").Value = "clorpromazina"
$ws.Range("C6").Value = "HALOPERIDOL"
$ws.Range("C7").Value = "clorpromazina"
$ws.Range("C8").Value = "codeina"
$ws.Range("C9").Value = "imipramina"
$ws.Range("C10").Value = "RISPERIDONA"
$ws.Range("C11").Value = "RISPERIDONA"
$ws.Range("C12").Value = "levomepromazina"
$ws.Range("C13").Value = "LIDOCAÍNA"
$ws.Range("C14").Value = "NITRATO DE CERIO;SULFADIAZINA"
$ws.Range("C15").Value = "COLAGENASE"

# Column D
$ws.Range("D2").Value = "5000ui"
$ws.Range("D3").Value = "25mg"
$ws.Range("D4").Value = "5mg"
$ws.Range("D5").Value = "40mg/ml"
$ws.Range("D6").Value = "2mg/ml"
$ws.Range("D7").Value = "25mg"
$ws.Range("D8").Value = "30mg"
$ws.Range("D9").Value = "25mg"
$ws.Range("D10").Value = "3mg"
$ws.Range("D11").Value = "1mg"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4%"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2%"
$ws.Range("D13").ClearFormats()
$ws.Range("D14").Value = "Concentração não encontrada"
$ws.Range("D15").Value = "30g"

# Column E
$ws.Range("E2").Value = "Cristália"
$ws.Range("E3").Value = "Teuto"
$ws.Range("E4").Value = "Cellera"
$ws.Range("E5").Value = "Sanofi"
$ws.Range("E6").Value = "Cellera"
$ws.Range("E7").Value = "Sanofi"
$ws.Range("E8").Value = "Cristália"
$ws.Range("E9").Value = "Cristália"
$ws.Range("E10").Value = "Vitamedic"
$ws.Range("E11").Value = "Vitamedic"
$ws.Range("E12").Value = "Sanofi"
$ws.Range("E13").Value = "Cristália"
$ws.Range("E14").Value = "Cristália"
$ws.Range("E15").Value = "Cristália"

# Column F
$ws.Range("F2").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Range("F3").Value = "LABORATÓRIO TEUTO BRASILEIRO S/A"
$ws.Range("F4").Value = "CELLERA FARMACÊUTICA S.A."
$ws.Range("F5").Value = "SANOFI MEDLEY FARMACÊUTICA LTDA"
$ws.Range("F6").Value = "CELLERA FARMACÊUTICA S.A."
$ws.Range("F7").Value = "SANOFI MEDLEY FARMACÊUTICA LTDA"
$ws.Range("F8").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Range("F9").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Range("F10").Value = "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA"
$ws.Range("F11").Value = "VITAMEDIC INDUSTRIA FARMACEUTICA LTDA"
$ws.Range("F12").Value = "SANOFI MEDLEY FARMACÊUTICA LTDA"
$ws.Range("F13").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Range("F14").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
$ws.Range("F15").Value = "CRISTÁLIA PRODUTOS QUÍMICOS FARMACÊUTICOS LTDA"
